$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for KWK Mysłowice-Wesoła Coal Mine, Poland, M1472, version ''Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on ' + $newStamp + ')''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

for ($row = 2; $row -le 11; $row++) {
    $cell = $wsData.Cells.Item($row, 19)
    $cell.Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
